$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (column G) values computed for rows 2-16 (header is row 1)
$newValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 4
    8  = 2
    9  = 0
    10 = 2
    11 = 3
    12 = 0
    13 = 2
    14 = 2
    15 = 3
    16 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
